$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "78÷6=13, 0"
$t.Cell(1, 2).Range.Text = "25÷6=4, 1"
$t.Cell(1, 3).Range.Text = "67÷3=22, 1"
$t.Cell(1, 4).Range.Text = "47÷2=23, 1"
$t.Cell(1, 5).Range.Text = "29÷6=4, 5"
$t.Cell(5, 1).Range.Text = "15÷5=3, 0"
$t.Cell(5, 2).Range.Text = "87÷4=21, 3"
$t.Cell(5, 3).Range.Text = "88÷7=12, 4"
$t.Cell(5, 4).Range.Text = "14÷6=2, 2"
$t.Cell(5, 5).Range.Text = "19÷4=4, 3"
$t.Cell(9, 1).Range.Text = "30÷6=5, 0"
$t.Cell(9, 2).Range.Text = "95÷3=31, 2"
$t.Cell(9, 3).Range.Text = "19÷9=2, 1"
$t.Cell(9, 4).Range.Text = "40÷7=5, 5"
$t.Cell(9, 5).Range.Text = "52÷5=10, 2"
$t.Cell(13, 1).Range.Text = "14÷5=2, 4"
$t.Cell(13, 2).Range.Text = "83÷6=13, 5"
$t.Cell(13, 3).Range.Text = "48÷8=6, 0"
$t.Cell(13, 4).Range.Text = "70÷9=7, 7"
$t.Cell(13, 5).Range.Text = "37÷3=12, 1"
$t.Cell(17, 1).Range.Text = "16÷2=8, 0"
$t.Cell(17, 2).Range.Text = "81÷6=13, 3"
$t.Cell(17, 3).Range.Text = "80÷4=20, 0"
$t.Cell(17, 4).Range.Text = "25÷3=8, 1"
$t.Cell(17, 5).Range.Text = "92÷6=15, 2"
